# sign-barcode-header.pptx edit
#  1. Bump the footer "date last updated" placeholder text on the slide
#     master and every slide layout from 6/10/2020 -> 6/16/2020.
#  2. Re-tint the two number-badge rectangles on slides 1-8 (a palette
#     refresh of the existing 8 slides).
#  3. Append 4 new number-badge slides (09-12), cloned from the slide 8
#     template, each with its own accent color.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Footer date placeholder: master + all custom layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($holder) {
    for ($i = 1; $i -le $holder.Shapes.Count; $i++) {
        $sh = $holder.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "6/16/2020"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------
# 2. Recolor the existing 8 slides' badge rectangles (shapes 1 & 2 on
#    every slide are "Rectangle 5" / "Rectangle 6").
# ---------------------------------------------------------------------
$existingColors = @{
    1 = 1241599   # FFF112 (was FFFF00)
    2 = 3617516   # EC3237 (was E24234)
    3 = 15773440  # 00AFF0 (was 203578)
    4 = 2714610   # F26B29 (was FCA457)
    5 = 5211904   # 00874F (was 346744)
    6 = 9381611   # EB268F (was FEB9CF)
    7 = 9781310   # 3E4095 (was 60539D)
    8 = 11910695  # 27BEB5 (was schemeClr tx2)
}

foreach ($idx in $existingColors.Keys) {
    $s = $p.Slides.Item($idx)
    $s.Shapes.Item(1).Fill.ForeColor.RGB = $existingColors[$idx]
    $s.Shapes.Item(2).Fill.ForeColor.RGB = $existingColors[$idx]
}

# ---------------------------------------------------------------------
# 3. Append 4 new slides (09, 10, 11, 12) cloned from slide 8's template.
# ---------------------------------------------------------------------
$newSlidesInfo = @(
    @{ Text = "09"; Rgb = 9130409 },  # A9518B
    @{ Text = "10"; Rgb = 6447200 },  # 606062
    @{ Text = "11"; Rgb = 1748724 },  # F4AE1A
    @{ Text = "12"; Rgb = 4640681 }   # A9CF46
)

foreach ($info in $newSlidesInfo) {
    $last = $p.Slides.Item($p.Slides.Count)
    $last.Duplicate()
    $new = $p.Slides.Item($p.Slides.Count)

    $new.Shapes.Item(1).Fill.ForeColor.RGB = $info.Rgb
    $new.Shapes.Item(2).Fill.ForeColor.RGB = $info.Rgb
    $new.Shapes.Item(3).TextFrame.TextRange.Text = $info.Text
}
